$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '71.739.19'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.72%  '
$ws.Range('E2').Style = 'Normal'

# Row 3 - Ethereum
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.981.66'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.68%  '
$ws.Range('E3').Style = 'Normal'

# Row 4 - TetherUSD
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('E4').Style = 'Normal'

# Row 5 - BNB
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '527.95'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('E5').Style = 'Normal'

# Row 6 - Solana
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '149.94'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.57%  '
$ws.Range('E6').Style = 'Normal'

# Row 7 - XRP
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.691'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +10.74%  '
$ws.Range('E7').Style = 'Normal'

# Row 8 - USDC
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E8').Style = 'Normal'

# Row 9 - Cardano
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.49%  '
$ws.Range('E9').Style = 'Normal'

# Row 10 - Dogecoin
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.53%  '
$ws.Range('E10').Style = 'Normal'

# Row 11 - ShibaInu
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000322'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -6.26%  '
$ws.Range('E11').Style = 'Normal'

# Row 12 - Avalanche
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.06'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +7.36%  '
$ws.Range('E12').Style = 'Normal'

# Row 13 - Polkadot
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.52'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.39%  '
$ws.Range('E13').Style = 'Normal'

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.618.03'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.44%  '
$ws.Range('E14').Style = 'Normal'

# Row 15 - WrappedEther
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.982.01'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.01%  '
$ws.Range('E15').Style = 'Normal'

# Row 16 - Uniswap
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.82'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -3.81%  '
$ws.Range('E16').Style = 'Normal'

# Row 17 - Chainlink
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '20.28'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -5.46%  '
$ws.Range('E17').Style = 'Normal'

# Row 18 - TRON
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.11%  '
$ws.Range('E18').Style = 'Normal'

# Row 19 - Polygon
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -4.58%  '
$ws.Range('E19').Style = 'Normal'

# Row 20 - WrappedBTC
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.545.47'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.48%  '
$ws.Range('E20').Style = 'Normal'

# Row 21 - BitcoinCash
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '423.63'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -4.04%  '
$ws.Range('E21').Style = 'Normal'

# Row 22 - Litecoin
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '97.35'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +4.63%  '
$ws.Range('E22').Style = 'Normal'

# Row 23 - ImmutableX
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.46'
$ws.Range('D23').Style = 'Normal'

# Row 24 - PancakeSwap
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.13'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.47%  '
$ws.Range('E24').Style = 'Normal'

# Row 25 - InternetComputer(DFINITY)
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '14.12'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.69%  '
$ws.Range('E25').Style = 'Normal'

# Row 26 - RenderToken
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -11.78%  '
$ws.Range('E26').Style = 'Normal'

# Row 27 - Filecoin
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.80%  '
$ws.Range('E27').Style = 'Normal'

# Row 28 - LEO
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.04%  '
$ws.Range('E28').Style = 'Normal'

# Row 29 - EthereumClassic
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '36.35'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.79%  '
$ws.Range('E29').Style = 'Normal'

# Row 30 - Toncoin
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.55'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +22.81%  '
$ws.Range('E30').Style = 'Normal'

# Row 31 - Cosmos
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.18'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -3.46%  '
$ws.Range('E31').Style = 'Normal'

# Row 32 - Bittensor
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '675.09'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.99%  '
$ws.Range('E32').Style = 'Normal'

# Row 33 - Hedera
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.59%  '
$ws.Range('E33').Style = 'Normal'

# Row 34 - NEARProtocol
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.85'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.25%  '
$ws.Range('E34').Style = 'Normal'

# Row 35 - OKB
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '64.82'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.56%  '
$ws.Range('E35').Style = 'Normal'

# Row 36 - InjectiveProtocol
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '42.03'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.38%  '
$ws.Range('E36').Style = 'Normal'

# Row 37 - TheGraph
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.421'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -5.16%  '
$ws.Range('E37').Style = 'Normal'

# Row 38 - PEPE
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0828'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -8.60%  '
$ws.Range('E38').Style = 'Normal'

# Row 39 - Kaspa
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.22%  '
$ws.Range('E39').Style = 'Normal'

# Row 40 - ThetaToken
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.43'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.86%  '
$ws.Range('E40').Style = 'Normal'

# Row 41 - Dai
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.16%  '
$ws.Range('E41').Style = 'Normal'

# Row 42 - WEMIXToken->FirstDigitalUSD (swap)
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('E42').Style = 'Normal'

# Row 43 - FirstDigitalUSD->WEMIXToken (swap)
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.28'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +5.35%  '
$ws.Range('E43').Style = 'Normal'

# Row 44 - VeChain
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.81%  '
$ws.Range('E44').Style = 'Normal'

# Row 45 - Stellar
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.55%  '
$ws.Range('E45').Style = 'Normal'

# Row 46 - THORChain
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.44%  '
$ws.Range('E46').Style = 'Normal'

# Row 47 - Fetch.AI
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.55'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -12.35%  '
$ws.Range('E47').Style = 'Normal'

# Row 48 - ApeXProtocol
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -5.36%  '
$ws.Range('E48').Style = 'Normal'

# Row 49 - Stacks
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -7.57%  '
$ws.Range('E49').Style = 'Normal'

# Row 50 - FLOKI
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000266'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -5.72%  '
$ws.Range('E50').Style = 'Normal'

# Row 51 - LidoDAOToken
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.24'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -5.00%  '
$ws.Range('E51').Style = 'Normal'
